$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 54, shifting existing rows 54-162 down to 55-163.
$ws.Rows.Item(54).Insert()

# Populate the newly inserted row 54 with the new data point.
$ws.Range("A54").Value = 11
$ws.Range("B54").Value = "Vega Monumental Concepción"
$ws.Range("C54").Value = "Bíobío"
$ws.Range("D54").Value = 44868
$ws.Range("E54").Value = 8
$ws.Range("F54").Value = 100112043
$ws.Range("G54").Value = "Pepino ensalada"
$ws.Range("H54").Value = "Sin especificar"
$ws.Range("I54").Value = "Primera"
$ws.Range("J54").Value = 100
$ws.Range("K54").Value = 17000
$ws.Range("L54").Value = 18000
$ws.Range("M54").Value = 17500
$ws.Range("N54").Value = "$/caja 60 unidades"
$ws.Range("O54").Value = "Región de Arica y Parinacota"
$ws.Range("P54").Value = 292
$ws.Range("Q54").Value = 60
$ws.Range("R54").Value = "Hortaliza"
